$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 previously held " - 16.06.2019"; simplify it to just the date.
$ws.Range("A2").Value = "16.06.2019"

# Add the two new task rows (7 and 8).
$ws.Range("A7").Value = "24.06.2019 - `n26.06.2019"
$ws.Range("B7").Value = "Working out how to connect backend with a `nminimal frontend app with proxy"
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = "26.06.2019 - `n27.06.2019"
$ws.Range("B8").Value = "Researching and testing Material-ui styles, tables and poppers and other related components for presenting lists of cards and their respective images"
$ws.Range("C8").Value = 6

# Match formatting of the existing wrapped/styled rows (e.g. row 5/row 4).
$ws.Range("A7:B8").WrapText = $true
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 45

# Move the active selection to B8, matching the final cursor position.
$ws.Range("B8").Select()
